$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 785
$ws.Range("I2").Value = 785
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 785
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -672

$ws.Range("H15").Value = 3309.1226
$ws.Range("I15").Value = 3309.1226
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 9927.3678
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -9758.3678

$ws.Range("H40").Value = 6281.8945
$ws.Range("I40").Value = 4449.1
$ws.Range("J40").Value = 8318.333000000001
$ws.Range("K40").Value = 4449.1
$ws.Range("L40").Value = 8318.333000000001
$ws.Range("M40").Value = -4274.1

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""

$ws.Range("H116").Value = 7066.067
$ws.Range("I116").Value = 6499.3
$ws.Range("J116").Value = 8199.6
$ws.Range("K116").Value = 6499.3
$ws.Range("L116").Value = 8199.6
$ws.Range("M116").Value = -3057.3
$ws.Range("N116").Value = -15083.6

$ws.Range("H125").Value = 432
$ws.Range("I125").Value = 432
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 3888
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -1428

$ws.Range("H137").Value = 2196.8293
$ws.Range("I137").Value = 2151.75
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 6455.25
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -3905.25

$ws.Range("H138").Value = 3620
$ws.Range("I138").Value = 3717.3333
$ws.Range("J138").Value = 3555.111
$ws.Range("K138").Value = 11151.9999
$ws.Range("L138").Value = 10665.333
$ws.Range("M138").Value = -6011.999899999999
$ws.Range("N138").Value = -20945.333

$ws.Range("H141").Value = 2189.0527
$ws.Range("I141").Value = 2205.6
$ws.Range("J141").Value = 2127
$ws.Range("K141").Value = 6616.799999999999
$ws.Range("L141").Value = 6381
$ws.Range("M141").Value = -1436.799999999999
$ws.Range("N141").Value = -16741

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3318.8235
$ws.Range("I61").Value = 1887.1428
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 1887.1428
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -1675.1428

$ws.Range("H63").Value = 2185.75
$ws.Range("I63").Value = 2185.75
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2185.75
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1499.75

$ws.Range("H66").Value = 2185.75
$ws.Range("I66").Value = 2185.75
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10928.75
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7496.75

$ws.Range("H94").Value = 50000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 50000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802

$ws.Range("H132").Value = 2441.5
$ws.Range("I132").Value = 2199
$ws.Range("J132").Value = 3249.8333
$ws.Range("K132").Value = 6597
$ws.Range("L132").Value = 9749.499899999999
$ws.Range("M132").Value = -4067

$ws.Range("H136").Value = 3318.8235
$ws.Range("I136").Value = 1887.1428
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 5661.428400000001
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -3111.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 99999
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 99999
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 99999
$ws.Range("N9").Value = -100335

$ws.Range("H86").Value = 2939.9375
$ws.Range("I86").Value = 2153.5
$ws.Range("J86").Value = 4250.6665
$ws.Range("K86").Value = 2153.5
$ws.Range("L86").Value = 4250.6665
$ws.Range("M86").Value = -1030.5

$ws.Range("H89").Value = 2939.9375
$ws.Range("I89").Value = 2153.5
$ws.Range("J89").Value = 4250.6665
$ws.Range("K89").Value = 10767.5
$ws.Range("L89").Value = 21253.3325
$ws.Range("M89").Value = -5151.5

$ws.Range("H134").Value = 6956.3335
$ws.Range("I134").Value = 2833.348
$ws.Range("J134").Value = 14250.846
$ws.Range("K134").Value = 8500.044
$ws.Range("L134").Value = 42752.538
$ws.Range("M134").Value = -5965.044
$ws.Range("N134").Value = -47822.538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 69294
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 69294
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 69294
$ws.Range("N68").Value = -70792

$ws.Range("H69").Value = 23591
$ws.Range("I69").Value = 12182
$ws.Range("J69").Value = 35000
$ws.Range("K69").Value = 12182
$ws.Range("L69").Value = 35000
$ws.Range("M69").Value = -11433
$ws.Range("N69").Value = -36498

$ws.Range("H71").Value = 69294
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 69294
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 207882
$ws.Range("N71").Value = -215370

$ws.Range("H72").Value = 23591
$ws.Range("I72").Value = 12182
$ws.Range("J72").Value = 35000
$ws.Range("K72").Value = 36546
$ws.Range("L72").Value = 105000
$ws.Range("M72").Value = -32802
$ws.Range("N72").Value = -112488

$ws.Range("H107").Value = 2942158.8
$ws.Range("I107").Value = 5000755.5
$ws.Range("J107").Value = 1306.4286
$ws.Range("K107").Value = 5000755.5
$ws.Range("L107").Value = 1306.4286
$ws.Range("M107").Value = -4998835.5
$ws.Range("N107").Value = -5146.4286

$ws.Range("H132").Value = 6539856
$ws.Range("I132").Value = 3719
$ws.Range("J132").Value = 37041828
$ws.Range("K132").Value = 11157
$ws.Range("L132").Value = 111125484
$ws.Range("M132").Value = -8627
$ws.Range("N132").Value = -111130544

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2775
$ws.Range("I35").Value = 2000
$ws.Range("J35").Value = 3550
$ws.Range("K35").Value = 6000
$ws.Range("L35").Value = 10650
$ws.Range("M35").Value = -5712
$ws.Range("N35").Value = -11226

$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 10000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622
$ws.Range("M69").Value = ""

$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 10000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112
$ws.Range("M72").Value = ""

$ws.Range("H82").Value = 11499.357
$ws.Range("I82").Value = 11500
$ws.Range("J82").Value = 11499.308
$ws.Range("K82").Value = 34500
$ws.Range("L82").Value = 34497.924
$ws.Range("M82").Value = -34094
$ws.Range("N82").Value = -35309.924

$ws.Range("H85").Value = 11499.357
$ws.Range("I85").Value = 11500
$ws.Range("J85").Value = 11499.308
$ws.Range("K85").Value = 34500
$ws.Range("L85").Value = 34497.924
$ws.Range("M85").Value = -33096
$ws.Range("N85").Value = -37305.924

$ws.Range("H122").Value = 1592.5714
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1592.5714
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 14333.1426
$ws.Range("N122").Value = -19233.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1139.8
$ws.Range("I3").Value = 1139.8
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1139.8
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1023.8

$ws.Range("H21").Value = 18999.666
$ws.Range("I21").Value = 6999
$ws.Range("J21").Value = 25000
$ws.Range("K21").Value = 6999
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = -6826
$ws.Range("N21").Value = -25346

$ws.Range("H29").Value = 2507
$ws.Range("I29").Value = 2507
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2507
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2217

$ws.Range("H30").Value = 18999.666
$ws.Range("I30").Value = 6999
$ws.Range("J30").Value = 25000
$ws.Range("K30").Value = 6999
$ws.Range("L30").Value = 25000
$ws.Range("M30").Value = -6894
$ws.Range("N30").Value = -25210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3633.5
$ws.Range("I46").Value = 2400.2222
$ws.Range("J46").Value = 7333.3335
$ws.Range("K46").Value = 2400.2222
$ws.Range("L46").Value = 7333.3335
$ws.Range("M46").Value = -2212.2222

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = ""

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = ""

$ws.Range("H100").Value = 7105.1
$ws.Range("I100").Value = 6300
$ws.Range("J100").Value = 7306.375
$ws.Range("K100").Value = 6300
$ws.Range("L100").Value = 7306.375
$ws.Range("M100").Value = -5759
$ws.Range("N100").Value = -8388.375

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""

$ws.Range("H136").Value = 3068
$ws.Range("I136").Value = 2858.8
$ws.Range("J136").Value = 3416.6667
$ws.Range("K136").Value = 8576.400000000001
$ws.Range("L136").Value = 10250.0001
$ws.Range("M136").Value = -6026.400000000001

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 13102.5
$ws.Range("I26").Value = 6956
$ws.Range("J26").Value = 19249
$ws.Range("K26").Value = 6956
$ws.Range("L26").Value = 19249
$ws.Range("M26").Value = -6663
$ws.Range("N26").Value = -19835

$ws.Range("H32").Value = 11500
$ws.Range("I32").Value = 11500
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 11500
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -11183

$ws.Range("H40").Value = 29341.666
$ws.Range("I40").Value = 29341.666
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 29341.666
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -29192.666
$ws.Range("N40").Value = ""

$ws.Range("H62").Value = 6556.5
$ws.Range("I62").Value = 5574.8
$ws.Range("J62").Value = 7538.2
$ws.Range("K62").Value = 5574.8
$ws.Range("L62").Value = 7538.2
$ws.Range("M62").Value = -4950.8

$ws.Range("H65").Value = 6556.5
$ws.Range("I65").Value = 5574.8
$ws.Range("J65").Value = 7538.2
$ws.Range("K65").Value = 27874
$ws.Range("L65").Value = 37691
$ws.Range("M65").Value = -24754

$ws.Range("H126").Value = 901
$ws.Range("I126").Value = 836
$ws.Range("J126").Value = 966
$ws.Range("K126").Value = 2508
$ws.Range("L126").Value = 2898
$ws.Range("M126").Value = -38
